$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: location/time updated, notes added
$ws.Range("C2").Value = "port willan"
$ws.Range("D2").Value = "2025-08-22 14:04:24"
$ws.Range("E2").Value = "big"

# Row 3: bird/number/location/time/notes updated
$ws.Range("A3").Value = "bellbird"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "3"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "sbhs"
$ws.Range("D3").Value = "2025-08-22 14:07:58"
$ws.Range("E3").Value = "breeding pair"

# Row 4: number/location/time updated (notes stays blank)
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "4"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "no"
$ws.Range("D4").Value = "2025-08-22 14:09:08"

# Rows 5-8 are removed entirely, shrinking the used range back to A1:E4
$ws.Range("A5:E8").Delete()
